$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.806.52"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "3.613.32"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'202.82"
$ws.Range("E5").Value = "  +3.79%  "
$ws.Range("D6").Value = "'599.20"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +6.69%  "
$ws.Range("D11").Value = "'53.99"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'0.0000304"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "'684.42"
$ws.Range("E14").Value = "  +15.06%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.181.72"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "70.869.26"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'19.26"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.84"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "3.607.88"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "'18.86"
$ws.Range("E22").Value = "  +6.02%  "
$ws.Range("D23").Value = "'110.32"
$ws.Range("E23").Value = "  +6.69%  "
$ws.Range("D24").Value = "'5.35"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("D25").Value = "'4.62"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'10.63"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").Value = "'6.02"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").Value = "'10.15"
$ws.Range("D30").Value = "'34.46"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("E31").Value = "  +5.84%  "
$ws.Range("D32").Value = "'7.19"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").Value = "'12.32"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'63.57"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "0.0₃0856"
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("D37").Value = "3.882.35"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'515.55"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  -4.90%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'37.04"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.60"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.385"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.139"
$ws.Range("E44").Value = "  +4.40%  "
$ws.Range("D45").Value = "'0.0470"
$ws.Range("E45").Value = "  +4.68%  "
$ws.Range("D46").Value = "'3.08"
$ws.Range("E46").Value = "  +9.64%  "
$ws.Range("D47").Value = "'3.41"
$ws.Range("E47").Value = "  +5.48%  "
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "'8.66"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'2.77"
$ws.Range("E51").Value = "  +68.94%  "
